$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("115:115").Insert()

$ws.Range("A115").Value = 10
$ws.Range("B115").Value = "Vega Modelo de Temuco"
$ws.Range("C115").Value = "La Araucanía"
$ws.Range("D115").Value = 44813
$ws.Range("E115").Value = 9
$ws.Range("F115").Value = "Fruta"
$ws.Range("G115").Value = 100104
$ws.Range("H115").Value = "Frutos de pepita"
$ws.Range("I115").Value = 100104001
$ws.Range("J115").Value = "Granada"
$ws.Range("K115").Value = "Wonderfull"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 55
$ws.Range("N115").Value = 14000
$ws.Range("O115").Value = 14000
$ws.Range("P115").Value = 14000
$ws.Range("Q115").Value = "$/bandeja 10 kilos granel"
$ws.Range("R115").Value = "Provincia de Limarí"
$ws.Range("S115").Value = 1400
$ws.Range("T115").Value = 10
